# The commit behind this diff ("Deployed fef2dc6 with MkDocs version: 1.3.1")
# is an automated docs-site deployment, not a manual spreadsheet edit.
#
# Diffing the OOXML shows no change to any cell value, formula, shared
# string, table, or merged range in any of the three sheets
# (算法复杂度 / 快速排序 / 归并排序). Every hunk is resave noise produced
# when the workbook was opened and saved by a different, genuine Excel
# build/platform than the one that wrote the previous copy:
#   - xl/workbook.xml: fileVersion/rupBuild, the x15ac:absPath (now a
#     macOS path), xr:revisionPtr documentId, and the bookViews window
#     geometry are all machine/session specific values stamped by the
#     authoring application on save.
#   - xl/styles.xml: the cellXfs table is reshuffled (two existing
#     <xf> entries move to the end of the list) and fonts gains
#     x14ac:knownFonts="1" - both are cosmetic artifacts of Excel's own
#     XML writer, not a formatting change (the two style records keep
#     their original numFmtId/font/fill/border/alignment content, only
#     their position - and therefore the s="" index referenced by the
#     few cells that use them - changes).
#   - xl/worksheets/sheet*.xml: sheetFormatPr gains baseColWidth="10"
#     and x14ac:dyDescent="0.2" (added to every row too) and
#     defaultRowHeight shifts 15.5 -> 16, which is the normal delta
#     between two Excel builds' default-font metrics; sheet3's column
#     widths move by the same kind of sub-pixel, metrics-driven amount.
#     None of this reflects a user action in the grid.
#
# There is therefore nothing for a script driving the Excel object
# model to author here - re-creating those values would mean faking a
# specific third-party machine's GUIDs/paths/screen geometry, which
# isn't something any Range/Worksheet/Workbook API call expresses, and
# the underlying cell data is already identical. The faithful action is
# simply to open the workbook (already done for us) and save it back
# out unchanged.

$wb = $excel.ActiveWorkbook
$wb.Save()
